$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16; existing rows 16.. shift down to 17..
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new weekly price entry
$ws.Cells.Item(16, 1).Value = 5
$ws.Cells.Item(16, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(16, 3).Value = "Maule"
$ws.Cells.Item(16, 4).Value = 45168
$ws.Cells.Item(16, 5).Value = 7
$ws.Cells.Item(16, 6).Value = 100112013
$ws.Cells.Item(16, 7).Value = "Alcachofa"
$ws.Cells.Item(16, 8).Value = "Madrigal"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 200
$ws.Cells.Item(16, 11).Value = 13000
$ws.Cells.Item(16, 12).Value = 13000
$ws.Cells.Item(16, 13).Value = 13000
$ws.Cells.Item(16, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(16, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(16, 16).Value = 325
$ws.Cells.Item(16, 17).Value = 40
$ws.Cells.Item(16, 18).Value = "Hortaliza"
